# Cropping-table workbook update:
#   - add "dz ind" computed-gap column between z1 ind (Q) and X (width) (old R)
#   - fill in the previously-empty "dx ind" (L) and "dy ind" (O) columns with
#     difference formulas (x1 ind - x0 ind, y1 ind - y0 ind)
#   - move the selection to R2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- dx ind (column L): L = K (x1 ind) - J (x0 ind) --------------------------
$ws.Range("L2").Formula = "=K2-J2"
$ws.Range("L3:L14").Formula = "=K3-J3"

# --- dy ind (column O): O = N (y1 ind) - M (y0 ind) --------------------------
$ws.Range("O2").Formula = "=N2-M2"
$ws.Range("O3:O14").Formula = "=N3-M3"

# --- new "dz ind" column, inserted before the old R (X (width)) column -------
$ws.Columns("R").Insert()
$ws.Range("R1").Value = "dz ind"

# --- selection ends up on R2 --------------------------------------------------
$ws.Range("R2").Select()
